$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Propagate the "last row" (heavy bottom-border) formatting that currently
#    sits on row 23 down to row 40 (the new last data row), BEFORE row 23's
#    own formatting gets reset to a normal interior row below.
# ---------------------------------------------------------------------------
$ws.Range("B23:J23").Copy()
$ws.Range("B40:J40").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Capture the two footer-note rows' formatting (currently rows 28/29)
#    before that area gets reused for table data, and relocate it to the
#    rows the footer will occupy after the table grows (45/46).
# ---------------------------------------------------------------------------
$ws.Range("B28:C28").Copy()
$ws.Range("B45:C45").PasteSpecial(-4122)
$ws.Range("H28:J28").Copy()
$ws.Range("H45:J45").PasteSpecial(-4122)

$ws.Range("B29:C29").Copy()
$ws.Range("B46:C46").PasteSpecial(-4122)
$ws.Range("H29:J29").Copy()
$ws.Range("H46:J46").PasteSpecial(-4122)

# The old footer rows were merged cells (B28:C28, H28:J28, B29:C29, H29:J29);
# those positions are about to hold normal per-cell table data, so break the
# merges before writing into them.
$ws.Range("B28:C28").UnMerge()
$ws.Range("B29:C29").UnMerge()
$ws.Range("H28:J28").UnMerge()
$ws.Range("H29:J29").UnMerge()

# The table-fill loop below only touches columns B-G, so explicitly blank out
# the old footer text that was sitting in H28/H29 (their merge-mates I/J were
# already empty).
$ws.Range("H28").Value = ""
$ws.Range("H29").Value = ""

# ---------------------------------------------------------------------------
# 3. Apply the normal data-row formatting (rows 16-22's style) to every row
#    that must look like a normal interior table row: the new rows 24-39,
#    plus row 23 (which stops being the last row for worker 1).
# ---------------------------------------------------------------------------
$ws.Range("B16:J16").Copy()
$ws.Range("B24:J39").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Fill the worker/period table (rows 16-40).
# ---------------------------------------------------------------------------
$workers = @(
    @("CC", "1128056533", "JAMES DAVID MARTINEZ RAMOS", 72318, 3883320, @("1712", "1711", "1710", "1709", "1708", "1707", "1706", "1705")),
    @("CC", "73571763", "FRANCISCO RAUL RADA MEJIA", 29509, 737717, @("1712", "1711", "1710", "1709", "1708", "1707", "1706", "1705")),
    @("CC", "73194932", "MIGUEL ANGEL BANQUEZ GONZALEZ", 29509, 737717, @("1712", "1711", "1710", "1709", "1708", "1707", "1706", "1705", "1704"))
)

$row = 16
foreach ($w in $workers) {
    $tipo = $w[0]
    $doc = $w[1]
    $nombre = $w[2]
    $salario = $w[3]
    $mora = $w[4]
    $periodos = $w[5]
    foreach ($p in $periodos) {
        $ws.Cells.Item($row, 2).Value = $tipo
        $ws.Cells.Item($row, 3).Value = $doc
        $ws.Cells.Item($row, 4).Value = $nombre
        $ws.Cells.Item($row, 5).Value = $p
        $ws.Cells.Item($row, 6).Value = $salario
        $ws.Cells.Item($row, 7).Value = $mora
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------------
# 5. Re-create the footer note rows (now 45/46) with their original text,
#    and re-merge the cells the way they were merged at 28/29.
# ---------------------------------------------------------------------------
$ws.Range("B45").Value = "___________________________________"
$ws.Range("H45").Value = "___________________________________"
$ws.Range("B46").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H46").Value = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B45:C45").Merge()
$ws.Range("H45:J45").Merge()
$ws.Range("B46:C46").Merge()
$ws.Range("H46:J46").Merge()

# ---------------------------------------------------------------------------
# 6. Update the summary block above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1080197
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 9

# ---------------------------------------------------------------------------
# 7. Column D widens to fit the new, longer worker names.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).AutoFit()
